$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Aydan's "Day 9 Finish" task cell: she is now covering for Manny
# instead of working on the final presentation.
$ws.Range("F6").Value = "Cover manny in his absence`nFinish unity sprite implementation for chicken run game `n"

# Re-autofit row 6 so the taller wrapped text doesn't leave a stray custom
# row height behind (content now fits the default row height again).
$ws.Rows(6).AutoFit()

# Widen the task columns so the fuller text is easier to read.
$ws.Columns("B").ColumnWidth = 109.42578125
$ws.Columns("D").ColumnWidth = 70.42578125
$ws.Columns("F").ColumnWidth = 70.42578125

# Move the active selection.
$ws.Range("F13").Select()
